$wb = $excel.ActiveWorkbook

# Update "展览" sheet (first sheet): F5 888 -> 890, F7 426 -> 427
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 890
$ws1.Range("F7").Value = 427

# Update "全部类型" sheet (fourth sheet): F5 888 -> 890, F7 426 -> 427
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 890
$ws4.Range("F7").Value = 427
